# Adds Harvard case classification results into the stats table.
#
# The previously-current "average_doctor"/"*_old" columns (and every other
# "*_old" app column affected by recomputed stats) are refreshed: the prior
# "average_doctor" figures become the new "average_doctor_old" baseline,
# while "average_doctor" itself (and the various "<app>_old" average /
# variance / std-dev columns) get the freshly recomputed values that include
# the new Harvard case.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "BP1" = "average_doctor_old"
    "BQ1" = "average_doctor"
    "AI4" = 0.288
    "AJ4" = 0.088
    "AK4" = 0.297
    "AU4" = 0.188
    "AV4" = 0.028
    "AW4" = 0.168
    "BA4" = 2.012
    "BB4" = 0.158
    "BC4" = 0.397
    "BG4" = 0.732
    "BH4" = 0.141
    "BI4" = 0.375
    "BM4" = 0.719
    "BN4" = 0.08
    "BO4" = 0.282
    "BP4" = 0.671
    "BQ4" = 0.71
    "E4" = 0.423
    "F4" = 0.072
    "G4" = 0.268
    "N4" = 0.432
    "O4" = 0.064
    "P4" = 0.254
    "Q4" = 0.023
    "R4" = 0.016
    "S4" = 0.128
    "W4" = 0.283
    "X4" = 0.109
    "Y4" = 0.33
    "AI5" = 0.309
    "AJ5" = 0.098
    "AK5" = 0.314
    "AU5" = 0.365
    "AV5" = 0.093
    "AW5" = 0.306
    "BA5" = 1.337
    "BB5" = 0.079
    "BC5" = 0.28
    "BG5" = 0.398
    "BH5" = 0.051
    "BI5" = 0.225
    "BM5" = 0.549
    "BN5" = 0.063
    "BO5" = 0.252
    "BP5" = 0.446
    "BQ5" = 0.456
    "E5" = 0.538
    "F5" = 0.085
    "G5" = 0.292
    "N5" = 0.733
    "O5" = 0.081
    "P5" = 0.285
    "Q5" = 0.015
    "R5" = 0.007
    "S5" = 0.082
    "W5" = 0.273
    "X5" = 0.109
    "Y5" = 0.33
    "AI6" = 0.298
    "AU6" = 0.248
    "BA6" = 1.599
    "BG6" = 0.516
    "BM6" = 0.623
    "BP6" = 0.533
    "BQ6" = 0.552
    "E6" = 0.474
    "N6" = 0.544
    "Q6" = 0.018
    "W6" = 0.278
    "AI7" = 0.305
    "AU7" = 0.307
    "BA7" = 1.429
    "BG7" = 0.438
    "BM7" = 0.576
    "BP7" = 0.476
    "BQ7" = 0.49
    "E7" = 0.51
    "N7" = 0.643
    "Q7" = 0.016
    "W7" = 0.275
    "AI8" = 0.33
    "AJ8" = 0.129
    "AK8" = 0.359
    "AU8" = 0.308
    "AV8" = 0.085
    "AW8" = 0.291
    "BA8" = 1.75
    "BB8" = 0.124
    "BC8" = 0.352
    "BG8" = 0.568
    "BH8" = 0.106
    "BI8" = 0.326
    "BM8" = 0.691
    "BN8" = 0.066
    "BO8" = 0.258
    "BP8" = 0.583
    "BQ8" = 0.605
    "E8" = 0.605
    "F8" = 0.112
    "G8" = 0.334
    "N8" = 0.773
    "O8" = 0.066
    "P8" = 0.256
    "Q8" = 0.017
    "S8" = 0.107
    "W8" = 0.301
    "Y8" = 0.347
    "AI9" = 0.255
    "AJ9" = 0.19
    "AK9" = 0.436
    "BA9" = 1.712
    "BB9" = 0.248
    "BC9" = 0.498
    "BG9" = 0.606
    "BH9" = 0.239
    "BI9" = 0.489
    "BM9" = 0.649
    "BN9" = 0.228
    "BO9" = 0.477
    "BP9" = 0.571
    "BQ9" = 0.588
    "E9" = 0.553
    "F9" = 0.247
    "G9" = 0.497
    "N9" = 0.681
    "O9" = 0.217
    "P9" = 0.466
    "W9" = 0.202
    "X9" = 0.161
    "Y9" = 0.402
    "AI10" = 0.362
    "AJ10" = 0.231
    "AK10" = 0.48
    "AU10" = 0.298
    "AV10" = 0.209
    "AW10" = 0.457
    "BA10" = 2.085
    "BB10" = 0.243
    "BC10" = 0.493
    "BG10" = 0.66
    "BH10" = 0.225
    "BI10" = 0.474
    "BM10" = 0.84
    "BN10" = 0.134
    "BO10" = 0.366
    "BP10" = 0.695
    "BQ10" = 0.726
    "E10" = 0.681
    "F10" = 0.217
    "G10" = 0.466
    "N10" = 0.872
    "O10" = 0.111
    "P10" = 0.334
    "W10" = 0.372
    "X10" = 0.234
    "Y10" = 0.483
    "AI11" = 0.394
    "AJ11" = 0.239
    "AK11" = 0.489
    "AU11" = 0.436
    "AV11" = 0.246
    "AW11" = 0.496
    "BA11" = 2.085
    "BB11" = 0.243
    "BC11" = 0.493
    "BG11" = 0.66
    "BH11" = 0.225
    "BI11" = 0.474
    "BM11" = 0.84
    "BN11" = 0.134
    "BO11" = 0.366
    "BP11" = 0.695
    "BQ11" = 0.728
    "E11" = 0.713
    "F11" = 0.205
    "G11" = 0.452
    "N11" = 0.894
    "O11" = 0.095
    "P11" = 0.308
    "W11" = 0.372
    "X11" = 0.234
    "Y11" = 0.483
    "AI12" = 1.703
    "AJ12" = 1.29
    "AK12" = 1.136
    "AU12" = 2.767
    "AV12" = 2.737
    "AW12" = 1.654
    "BB12" = 0.395
    "BC12" = 0.629
    "BG12" = 1.097
    "BH12" = 0.12
    "BI12" = 0.346
    "BM12" = 1.291
    "BN12" = 0.333
    "BO12" = 0.577
    "BP12" = 1.232
    "BQ12" = 1.257
    "E12" = 1.403
    "F12" = 0.748
    "G12" = 0.865
    "N12" = 1.465
    "O12" = 1.039
    "P12" = 1.02
    "W12" = 1.629
    "X12" = 0.576
    "Y12" = 0.759
    "AI13" = 1.28
    "AJ13" = 0.37
    "AK13" = 0.608
    "AU13" = 2.285
    "AV13" = 0.925
    "AW13" = 0.962
    "BA13" = 2.344
    "BB13" = 0.296
    "BC13" = 0.544
    "BG13" = 0.584
    "BH13" = 0.071
    "BI13" = 0.266
    "BM13" = 0.892
    "BN13" = 0.282
    "BO13" = 0.531
    "BP13" = 0.781
    "BQ13" = 0.725
    "E13" = 1.573
    "F13" = 0.652
    "G13" = 0.808
    "N13" = 2.068
    "O13" = 0.928
    "P13" = 0.964
    "W13" = 1.037
    "X13" = 0.193
    "Y13" = 0.439
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}